$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6638576779026217
$ws1.Range("C2").Value = 0.6140808344198174
$ws1.Range("D2").Value = 0.8820224719101124
$ws1.Range("E2").Value = 0.7240584166026134
$ws1.Range("F2").Value = 0.8112297623148467
$ws1.Range("G2").Value = 0.867464758801445
$ws1.Range("H2").Value = 0.6710782869727447
$ws1.Range("I2").Value = 471
$ws1.Range("J2").Value = 296
$ws1.Range("K2").Value = 238
$ws1.Range("L2").Value = 63

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.7906976744186046
$ws2.Range("C2").Value = 0.4456928838951311
$ws2.Range("D2").Value = 0.5700598802395209

$ws2.Range("B3").Value = 0.6140808344198174
$ws2.Range("C3").Value = 0.8820224719101124
$ws2.Range("D3").Value = 0.7240584166026134

$ws2.Range("B4").Value = 0.6638576779026217
$ws2.Range("C4").Value = 0.6638576779026217
$ws2.Range("D4").Value = 0.6638576779026217
$ws2.Range("E4").Value = 0.6638576779026217

$ws2.Range("B5").Value = 0.702389254419211
$ws2.Range("C5").Value = 0.6638576779026217
$ws2.Range("D5").Value = 0.6470591484210672

$ws2.Range("B6").Value = 0.7023892544192111
$ws2.Range("C6").Value = 0.6638576779026217
$ws2.Range("D6").Value = 0.647059148421067

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 238
$ws3.Range("C2").Value = 296

$ws3.Range("B3").Value = 63
$ws3.Range("C3").Value = 471
